$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 07:35"

# Update Hungria (row 69) stats
$ws.Range("B69").Value = 3535
$ws.Range("C69").Value = 26
$ws.Range("D69").Value = 1400
$ws.Range("E69").Value = 1673
$ws.Range("G69").Value = 11
$ws.Range("H69").Value = 462

# Update Uzbekistan (row 75) stats
$ws.Range("B75").Value = 2762
$ws.Range("C75").Value = 9
$ws.Range("E75").Value = 503

# Row 81 previously showed "Croacia"; it now becomes "Bulgaria" with updated stats
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 2235
$ws.Range("C81").Value = 24
$ws.Range("D81").Value = 612
$ws.Range("E81").Value = 1513
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 110

# Row 82 previously showed "Bulgaria"; it now becomes "Croacia" keeping the old Croacia stats
$ws.Range("A82").Value = "Croacia"
$ws.Range("B82").Value = 2226
$ws.Range("D82").Value = 1936
$ws.Range("E82").Value = 195
$ws.Range("H82").Value = 95
